$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    9  = -10.331
    18 = -12.89079999999999
    20 = -11.50000000000001
    27 = -12.548
    35 = -11.7556
    69 = -11.3874
    76 = -12.22740000000001
    78 = -11.48010000000001
    82 = -11.86289999999999
    83 = -14.1146
    93 = -10.9819
}

foreach ($row in $updates.Keys) {
    $ws.Range("C$row").Value = $updates[$row]
}
